$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.286832544864788, 10.34677158129881, 0.1494219747398047, 10.19245300693656, 0, 23.97547910783996)
    3 = @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 6.189590430959694)
    4 = @(0.6606524410359556, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 2.960089034096801)
    5 = @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 3.754798637575387)
    6 = @(3.286832544864788, 1.655778082260271, 3.537761648806719, 10.19245300693656, 0, 18.67282528286833)
    7 = @(0.6606524410359556, 0.306821227259698, 0.7527432677738641, 0.4942365360607697, 0, 2.214453472130288)
    8 = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 0, 7.143138311642302)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
